# Added Crdc login backup codes
#
# The sheet previously held 3 placeholder backup codes at A2:A4 and the
# "real" codes starting at A11. This edit promotes the first three real
# codes (formerly at A11:A13) up into A2:A4, and removes the now-redundant
# A11:A13 rows (the remaining codes at A14:A17 are untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Promote the codes that used to live at A11:A13 into A2:A4.
$ws.Range("A2").Value = $ws.Range("A11").Value()
$ws.Range("A3").Value = $ws.Range("A12").Value()
$ws.Range("A4").Value = $ws.Range("A13").Value()

# Remove the now-duplicated rows 11-13 (clearing makes the rows disappear
# from sheetData entirely, rather than shifting rows 14:17 upward).
$ws.Range("A11:A13").ClearContents()

# Match the author's final selection.
$ws.Range("A10").Select()
